$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at K, pushing "Web Query" (old K) and everything after it one column right.
$ws.Columns("K:K").Insert()

# Insert a new blank column at M (after "Web Query", now in L), pushing "Number of Proposed Services"
# (now in M) and everything after it one column right.
$ws.Columns("M:M").Insert()

# Fill in the new header cells.
$ws.Range("K1").Value2 = "Web Search Mode"
$ws.Range("M1").Value2 = "Attack Type"

# Fill in the new data columns for each data row (rows 2-21).
$attackTypes = @{
    2  = "Remote Code Execution"
    3  = "Remote Code Execution"
    4  = "Local File Inclusion (LFI) leading to Remote Code Execution (RCE)"
    5  = "Remote Code Execution"
    6  = "remote code execution"
    7  = "Unauthorized access to sensitive files"
    8  = "Remote Code Execution"
    9  = "Information Disclosure"
    10 = "Path Traversal"
    11 = "Path Traversal and Remote Code Execution"
    12 = "Path Traversal and Remote Code Execution"
    13 = "Path Traversal"
    14 = "Remote Code Execution"
    15 = "Code Injection, Remote Code Execution"
    16 = "Remote Code Execution"
    17 = "Remote Code Execution"
    18 = "Remote Code Execution"
    19 = "Information Disclosure"
    20 = "unauthenticated remote code execution"
    21 = "Remote Code Execution"
}

for ($row = 2; $row -le 21; $row++) {
    $ws.Range("K$row").Value2 = "custom_no_tool"
    $ws.Range("M$row").Value2 = $attackTypes[$row]
}
